# Auto-generated Excel COM-interop script applying scheduled-runner market data updates
# to the Phoenix_Profits workbook (regenerates currentAveragePrice / Leve profit columns).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (87 value updates) ----
$ws = $wb.Worksheets.Item("ALC")

$edits = @(
    @("H9", 20419.8),
    @("I9", 49.5),
    @("J9", 34000),
    @("K9", 49.5),
    @("L9", 34000),
    @("M9", 119.5),
    @("N9", -34338),
    @("H12", 254.8),
    @("I12", 225),
    @("K12", 225),
    @("M12", -55),
    @("H19", 3903.4517),
    @("I19", 605),
    @("J19", 5474.143),
    @("K19", 605),
    @("L19", 5474.143),
    @("M19", -430),
    @("N19", -5824.143),
    @("H28", 1869),
    @("I28", 1710.8889),
    @("K28", 1710.8889),
    @("M28", -1225.8889),
    @("H80", 329.9),
    @("I80", 163),
    @("K80", 489),
    @("M80", 509),
    @("H83", 329.9),
    @("I83", 163),
    @("K83", 1467),
    @("M83", 3525),
    @("H92", 20219.96),
    @("I92", 26435.947),
    @("J92", 536),
    @("K92", 26435.947),
    @("L92", 536),
    @("M92", -25187.947),
    @("N92", -3032),
    @("H98", 42086.293),
    @("I98", 51997.26),
    @("J98", 4424.6),
    @("K98", 51997.26),
    @("L98", 4424.6),
    @("M98", -50499.26),
    @("N98", -7420.6),
    @("H103", 880.6539),
    @("I103", 443.07144),
    @("J103", 1391.1666),
    @("K103", 1329.21432),
    @("L103", 4173.4998),
    @("M103", -743.21432),
    @("N103", -5345.4998),
    @("H104", 1826.8334),
    @("I104", 1826.8334),
    @("K104", 5480.5002),
    @("M104", -3733.5002),
    @("H113", 2173.8),
    @("I113", 1957.3334),
    @("J113", 2498.5),
    @("K113", 1957.3334),
    @("L113", 2498.5),
    @("M113", 1296.6666),
    @("N113", -9006.5),
    @("H116", 4962.7334),
    @("I116", 4236.4194),
    @("J116", 6571),
    @("K116", 4236.4194),
    @("L116", 6571),
    @("M116", -794.4193999999998),
    @("N116", -13455),
    @("H122", 42086.293),
    @("I122", 51997.26),
    @("J122", 4424.6),
    @("K122", 155991.78),
    @("L122", 13273.8),
    @("M122", -153541.78),
    @("N122", -18173.8),
    @("H132", 3609),
    @("I132", 3777.45),
    @("K132", 11332.35),
    @("M132", -8802.349999999999),
    @("H138", 2956.0652),
    @("I138", 2031.1111),
    @("J138", 3550.6785),
    @("K138", 6093.3333),
    @("L138", 10652.0355),
    @("M138", -953.3333000000002),
    @("N138", -20932.0355)
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}


# ---- Sheet: ARM (22 value updates) ----
$ws = $wb.Worksheets.Item("ARM")

$edits = @(
    @("H32", 3302.0527),
    @("I32", 3030.6826),
    @("K32", 3030.6826),
    @("M32", -2743.6826),
    @("H56", 361662.34),
    @("J56", 50000),
    @("L56", 50000),
    @("N56", -51484),
    @("H97", 30303794),
    @("I97", 417.77274),
    @("J97", 90910540),
    @("K97", 417.77274),
    @("L97", 90910540),
    @("M97", 78.22726),
    @("N97", -90911532),
    @("H132", 2166.8975),
    @("I132", 2085.8823),
    @("J132", 2717.8),
    @("K132", 6257.646900000001),
    @("L132", 8153.400000000001),
    @("M132", -3727.646900000001),
    @("N132", -13213.4)
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}


# ---- Sheet: BSM (8 value updates) ----
$ws = $wb.Worksheets.Item("BSM")

$edits = @(
    @("H94", 82587.45),
    @("I94", 766.125),
    @("K94", 766.125),
    @("M94", -315.125),
    @("H107", 3557.3215),
    @("I107", 3654.077),
    @("K107", 3654.077),
    @("M107", -1734.077)
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}


# ---- Sheet: CRP (26 value updates) ----
$ws = $wb.Worksheets.Item("CRP")

$edits = @(
    @("H31", 2149.7778),
    @("I31", 1869.3334),
    @("K31", 1869.3334),
    @("M31", -1574.3334),
    @("H34", 2149.7778),
    @("I34", 1869.3334),
    @("K34", 1869.3334),
    @("M34", -1667.3334),
    @("H58", 1922.5682),
    @("I58", 1540.375),
    @("J58", 2941.75),
    @("K58", 1540.375),
    @("L58", 2941.75),
    @("M58", -1337.375),
    @("N58", -3347.75),
    @("H94", 1713.1578),
    @("I94", 1989.1818),
    @("K94", 1989.1818),
    @("M94", -1538.1818),
    @("H136", 1922.5682),
    @("I136", 1540.375),
    @("J136", 2941.75),
    @("K136", 4621.125),
    @("L136", 8825.25),
    @("M136", -2071.125),
    @("N136", -13925.25)
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}


# ---- Sheet: CUL (12 value updates) ----
$ws = $wb.Worksheets.Item("CUL")

$edits = @(
    @("H18", 1534.0714),
    @("I18", 281.625),
    @("K18", 844.875),
    @("M18", -675.875),
    @("H50", 978),
    @("I50", 974.2),
    @("K50", 2922.6),
    @("M50", -2441.6),
    @("H53", 978),
    @("I53", 974.2),
    @("K53", 2922.6),
    @("M53", -2441.6)
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}


# ---- Sheet: GSM (22 value updates, 1 cleared cell(s)) ----
$ws = $wb.Worksheets.Item("GSM")

$edits = @(
    @("H59", 7800),
    @("J59", 7800),
    @("L59", 7800),
    @("N59", -8966),
    @("H96", 12999.5),
    @("J96", 12999.5),
    @("L96", 12999.5),
    @("N96", -18491.5),
    @("H102", 20450.133),
    @("I102", 22192.285),
    @("K102", 22192.285),
    @("M102", -20570.285),
    @("H132", 3037.52),
    @("I132", 2914.0833),
    @("J132", 6000),
    @("K132", 8742.249899999999),
    @("L132", 18000),
    @("M132", -6212.249899999999),
    @("N132", -23060),
    @("H133", 0),
    @("J133", 0),
    @("L133", 0)
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

$ws.Range("N133").ClearContents()

# ---- Sheet: LTW (35 value updates, 1 cleared cell(s)) ----
$ws = $wb.Worksheets.Item("LTW")

$edits = @(
    @("H16", 942.98114),
    @("I16", 893.3022999999999),
    @("J16", 1156.6),
    @("K16", 893.3022999999999),
    @("L16", 1156.6),
    @("M16", -723.3022999999999),
    @("N16", -1496.6),
    @("H46", 7418.2),
    @("I46", 0),
    @("J46", 7418.2),
    @("K46", 0),
    @("L46", 7418.2),
    @("N46", -7794.2),
    @("H82", 1209.4615),
    @("I82", 974.1),
    @("K82", 974.1),
    @("M82", -613.1),
    @("H85", 1209.4615),
    @("I85", 974.1),
    @("K85", 974.1),
    @("M85", 273.9),
    @("H100", 2249.6956),
    @("I100", 1793.7894),
    @("J100", 4415.25),
    @("K100", 1793.7894),
    @("L100", 4415.25),
    @("M100", -1252.7894),
    @("N100", -5497.25),
    @("H132", 3098.75),
    @("I132", 2210.16),
    @("J132", 5118.273),
    @("K132", 6630.48),
    @("L132", 15354.819),
    @("M132", -4100.48),
    @("N132", -20414.819)
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

$ws.Range("M46").ClearContents()

# ---- Sheet: WVR (23 value updates) ----
$ws = $wb.Worksheets.Item("WVR")

$edits = @(
    @("H100", 52632400),
    @("I100", 83334110),
    @("K100", 166668220),
    @("M100", -166667679),
    @("H107", 395.69232),
    @("I107", 293.05554),
    @("K107", 879.16662),
    @("M107", 1040.83338),
    @("H113", 1161.0741),
    @("I113", 1271.7368),
    @("K113", 3815.2104),
    @("M113", -1645.2104),
    @("H126", 45460310),
    @("I126", 50005664),
    @("K126", 150016992),
    @("M126", -150014522),
    @("H132", 2774.6128),
    @("I132", 2655.6206),
    @("J132", 4500),
    @("K132", 7966.861800000001),
    @("L132", 13500),
    @("M132", -5436.861800000001),
    @("N132", -18560)
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

